$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45171 -> 45172) for every data row (rows 2 through 344).
$ws.Range("C2:C344").Value = 45172
